$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")
Write-Host $ws.Name
